$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Черепанова"
$ws.Range("A2").Value = "Иванова"
$ws.Range("A1").Value = "Петров"
$ws.Range("A4").Value = "Мещеряков"
$ws.Range("A5").Value = "Голубев"
$ws.Range("A6").Value = "Орлов"

$ws.Columns.Item(1).ColumnWidth = 15.6

$ws.Range("A7").Select()
